# Populate the 'IIS errors' reference table on the active sheet,
# matching the layout of the finished workbook: a bold header row,
# an unsorted data entry pass, final A-Z sort by error code, and
# the column widths / row heights / page setup used for printing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header row ---------------------------------------------------
$ws.Range("A1").Value = 'error code'
$ws.Range("B1").Value = 'Description'

# --- data rows, entered in original (pre-sort) order ---------------
$ws.Range("A2").Value = 500
$ws.Range("B2").Value = 'Ошибка с кодом «500» - возникла в результате неверного синтаксиса файла .htaccess или наличие в нем неподдерживаемых директив, а так же вследствие неправильного обращения с CGI-скриптами.'
$ws.Range("A3").Value = 401
$ws.Range("B3").Value = 'Ошибка с кодом «401» - доступ запрещен списком управления доступов к ресурсу.'
$ws.Range("A4").Value = 503
$ws.Range("B4").Value = 'Ошибка с кодом «503» - Эта ошибка возникает, если не запустить пул приложений, связанный с веб-приложением. Чтобы устранить эту проблему, выполните следующие действия.'
$ws.Range("A5").Value = 404
$ws.Range("B5").Value = 'Ошибка с кодом «404» - возникла в результате неверных url-ссылок на файлы, расположенные в веб-приложении. Причины возникновения: 
в коде веб-приложения указаны ошибочные наименования (никогда ранее не существующих) файлов; 
указанные в логах файлы ранее существовали (использовались) в веб-приложении, но в результате модификации (в различных формах - развитие, оптимизация) веб-приложения функции, которые были описаны в файлах были модифицированы, а ссылки на устаревшие файлы из кода веб-приложения не было полностью удалены, отправлен разработчику на устранение ошибок.'
$ws.Range("A6").Value = 403
$ws.Range("B6").Value = 'Ошибка с кодом «403» - ошибка запрета доступа. Причины возникновения:
на сервере имеется неправильный индексный файл;
выставлены не верные права на папку в которой находится запрашиваемый файл, либо на какую-либо из ее родительских директорий;
файл загружен в неправильную папку.'
$ws.Range("A7").Value = 400
$ws.Range("B7").Value = 'Ошибка с кодом «400» - ошибка доступа к серверу шлюза TS. Причины возникновения: услуги, требуемые шлюзом TS, не запускаются; существуют проблемы с сервером NPS и Web-сервером (IIS).'
$ws.Range("A8").Value = 405
$ws.Range("B8").Value = 'Ошибка с кодом «405» - означает, что указанный клиентом метод нельзя применить к текущему ресурсу.'
$ws.Range("A9").Value = 406
$ws.Range("B9").Value = 'Ошибка с кодом «406» - означает, что браузер клиента не принимает тип MIME запрашиваемой страницы.'
$ws.Range("A10").Value = 101
$ws.Range("B10").Value = 'Ошибка с кодом «101» - означает, что сервер выполняет требование клиента и переключает протоколы в соответствии с указанием, данным в поле заголовка Upgrade. Сервер отправляет заголовок ответа Upgrade, указывая протокол, на который он переключился.'
$ws.Range("A11").Value = 503
$ws.Range("B11").Value = 'Ошибка с кодом «503» - эта ошибка возникает, если не запустить пул приложений, связанный с веб-приложением. Чтобы устранить эту проблему, выполните следующие действия.'
$ws.Range("A12").Value = 302
$ws.Range("B12").Value = 'Ошибка с кодом «302» - означает, что запрошенный документ временно доступен по другому URI, указанному в заголовке в поле Location. Этот код может быть использован, например, при управляемом сервером согласовании содержимого.'
$ws.Range("A13").Value = 304
$ws.Range("B13").Value = 'Ошибка с кодом «304» - означает, что сервер возвращает такой код, если клиент запросил документ методом GET, использовал заголовок If-Modified-Since или If-None-Match и документ не изменился с указанного момента. При этом сообщение сервера не должно содержать тела.'
$ws.Range("A14").Value = 305
$ws.Range("B14").Value = 'Ошибка с кодом «305» - означает, что запрос к запрашиваемому ресурсу должен осуществляться через прокси-сервер, URI которого указан в поле Location заголовка. Данный код ответа могут использовать только исходные HTTP-сервера (не прокси)'
$ws.Range("A15").Value = 301
$ws.Range("B15").Value = 'Ошибка с кодом «301» - означает, что запрошенный документ был окончательно перенесен на новый URI, указанный в поле Location заголовка. Некоторые клиенты некорректно ведут себя при обработке данного кода.'
$ws.Range("A16").Value = 206
$ws.Range("B16").Value = 'Код состояния «206» -  это ответ на запрос части документа. Это используется расширенными инструментами кэширования, когда пользовательский агент запрашивает только небольшую часть страницы, и возвращается только этот раздел.'

# --- alignment on the data body (rows 2:16), applied before the ---
# --- header is touched so the header gets its own plain-bold xf ---
$ws.Range("B2:B16").WrapText = $true
$ws.Range("A2:A16").HorizontalAlignment = -4108

# error 404 (pre-sort row 5) uses top-aligned wrapped text instead
$ws.Range("B5").VerticalAlignment = -4160
$ws.Range("B5").WrapText = $true

# --- header styling: bold first, then its own alignment -----------
$ws.Range("A1:B1").Font.Bold = $true
$ws.Range("B1").WrapText = $true
$ws.Range("A1").HorizontalAlignment = -4108

# --- sort the data body by error code, ascending ------------------
# (uses the persistent Sort object so the workbook keeps a sortState
# the same way Excel's Data > Sort dialog does)
$ws.Sort.SortFields.Clear()
[void]$ws.Sort.SortFields.Add($ws.Range("A1"))
$ws.Sort.SetRange($ws.Range("A1:B16"))
$ws.Sort.Header = 1
[void]$ws.Sort.Apply()

# --- column widths --------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 18.85546875
$ws.Columns.Item(2).ColumnWidth = 173.42578125

# --- row heights (final, post-sort positions) -----------------------
$ws.Rows.Item(2).RowHeight = 30
$ws.Rows.Item(3).RowHeight = 30
$ws.Rows.Item(4).RowHeight = 30
$ws.Rows.Item(5).RowHeight = 30
$ws.Rows.Item(6).RowHeight = 57.95
$ws.Rows.Item(7).RowHeight = 30
$ws.Rows.Item(8).RowHeight = 29.1
$ws.Rows.Item(10).RowHeight = 60
$ws.Rows.Item(11).RowHeight = 75
$ws.Rows.Item(14).RowHeight = 30

# --- selection + view --------------------------------------------
[void]$ws.Range("B6").Select()

# --- page setup ------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

